# "using lockdown for UK, Europe"
#
# 1. Update mitigation-date sources for several European countries + UK
#    (rows 2-12 on the "mitigation" sheet) to newer lockdown dates.
# 2. Record that the Imperial College report source (row 1 on
#    "mitigation sources") is tagged as "lockdown".
# 3. Drop the empty "recovery_statistics" sheet and rename
#    "death_statistics" to "rec_death_statistics", hiding its (now
#    redundant) recovery source row.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false | Out-Null

# --- 1. mitigation sheet: refreshed lockdown dates -------------------------
$ws1 = $wb.Worksheets("mitigation")

$ws1.Range("B2").Value2  = 43906   # Austria
$ws1.Range("B3").Value2  = 43908   # Belgium
$ws1.Range("B4").Value2  = 43908   # Denmark
$ws1.Range("B5").Value2  = 43907   # France
$ws1.Range("B6").Value2  = 43912   # Germany
$ws1.Range("B7").Value2  = 43901   # Italy
$ws1.Range("B8").Value2  = 43914   # Norway
$ws1.Range("B9").Value2  = 43904   # Spain
$ws1.Range("B10").Value2 = 43906   # Sweden
$ws1.Range("B11").Value2 = 43910   # Switzerland
$ws1.Range("B12").Value2 = 43914   # United Kingdom

$ws1.Range("D8").Select() | Out-Null

# --- 2. mitigation sources sheet: tag Europe/UK source as "lockdown" -------
$ws2 = $wb.Worksheets("mitigation sources")
$ws2.Range("B1").Value = "lockdown"

# --- 3. drop recovery_statistics, rename death_statistics -------------------
# (rename before deleting the other sheet: deleting a sheet shifts sheet
#  indices, and an already-fetched Worksheet reference becomes stale)
$ws4 = $wb.Worksheets("death_statistics")
$ws4.Rows(2).Hidden = $true
$ws4.Range("A6").Select() | Out-Null
$ws4.Name = "rec_death_statistics"

$wb.Worksheets("recovery_statistics").Delete() | Out-Null

# --- finally, make "mitigation sources" the active/visible tab -------------
$ws2.Activate() | Out-Null
$ws2.Range("B1").Select() | Out-Null
